# Update the "no convexos" experiment data (alpha_zero generator) with the
# new values produced by the generator re-run. All affected cells hold
# literal text (the source data is exported from Python as strings, even
# when the text looks numeric), so every write below is done through a
# temporary "@" (Text) number format to stop Excel's automatic number
# coercion, and the format is cleared again immediately afterwards so the
# cells end up with no explicit style (matching their original state).

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ------------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $ws.Range("A2") "5.0 - y_1"
Set-TextValue $ws.Range("B2") "-5.0"
Set-TextValue $ws.Range("D2") "0.35"
Set-TextValue $ws.Range("E2") "9.1"
Set-TextValue $ws.Range("F2") "8.0"

Set-TextValue $ws.Range("A3") "-5.0 + y_1"
Set-TextValue $ws.Range("B3") "1.0"
Set-TextValue $ws.Range("D3") "0.91"
Set-TextValue $ws.Range("E3") "7.800000000000001"
Set-TextValue $ws.Range("F3") "6.4"

Set-TextValue $ws.Range("A4") "3.0000000000000036 - 2x + y_1 + 4y_2"
Set-TextValue $ws.Range("B4") "-19.000000000000004"
Set-TextValue $ws.Range("D4") "0.78"
Set-TextValue $ws.Range("E4") "1.2"
Set-TextValue $ws.Range("F4") "4.0"

Set-TextValue $ws.Range("A5") "-70.1 + 8x + y_1"
Set-TextValue $ws.Range("B5") "21.799999999999997"
Set-TextValue $ws.Range("D5") "0.57"
Set-TextValue $ws.Range("E5") "6.0"
Set-TextValue $ws.Range("F5") "7.9"

Set-TextValue $ws.Range("A6") "-2.1999999999999993 - 2x - 2y_1"
Set-TextValue $ws.Range("B6") "-14.2"
Set-TextValue $ws.Range("D6") "0.0"
Set-TextValue $ws.Range("E6") "2.1"

# --- Punto_modificado ------------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")

Set-TextValue $ws.Range("A2") "8.1"
Set-TextValue $ws.Range("B2") "5.0"
Set-TextValue $ws.Range("C2") "2.05"

# --- Vector_bf ---------------------------------------------------------
# NOTE: "Vector_bf" and "Vector_BF" only differ by case, and
# Worksheets.Item(<name>) resolves case-insensitively here, so both names
# would otherwise hit the same sheet. Use the (1-based) sheet index instead
# to land on the correct tab: sheet 5 = Vector_bf, sheet 6 = Vector_BF.
$ws = $wb.Worksheets.Item(5)

Set-TextValue $ws.Range("A2") "-0.9100000000000001"
Set-TextValue $ws.Range("A3") "-3.12"

# --- Vector_BF ---------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

Set-TextValue $ws.Range("A2") "-40.4"
Set-TextValue $ws.Range("A3") "1.299999999999999"
Set-TextValue $ws.Range("A4") "-6.8"
